# Update ShipmentTrackNum (column C) / PackageTrackNum (column D) values
# for rows 2-16 on the CheetahProcessing "Login email" refresh (12 Jul 2022).
#
# New tracking numbers keyed by row. For rows where the original workbook
# also mirrored the value into column D (rows 5,6,7,13,14,15,16), set D too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "320018680212"
    3  = "320018680223"
    4  = "320018680256"
    5  = "320018680278"
    6  = "320018680315"
    7  = "320018680337"
    8  = "320018680360"
    9  = "320018680381"
    10 = "320018680418"
    11 = "320018680430"
    12 = "320018680473"
    13 = "320018680495"
    14 = "320018680521"
    15 = "320018680543"
    16 = "320018680576"
}

$mirrorToD = @(5, 6, 7, 13, 14, 15, 16)

foreach ($row in $values.Keys) {
    $val = $values[$row]

    $cCell = $ws.Range("C$row")
    $cCell.Value = "'" + $val
    $cCell.Style = "Normal"

    if ($mirrorToD -contains $row) {
        $dCell = $ws.Range("D$row")
        $dCell.Value = "'" + $val
        $dCell.Style = "Normal"
    }
}
